$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.382.16"
$ws.Range("E2").Value = "  -7.81%  "
$ws.Range("D3").Value = "3.675.37"
$ws.Range("E3").Value = "  -7.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.93"
$ws.Range("E5").Value = "  -8.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.51"
$ws.Range("E6").Value = "  +5.04%  "
$ws.Range("D7").Value = "3.665.86"
$ws.Range("E7").Value = "  -7.44%  "
$ws.Range("E8").Value = "  -7.41%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -6.34%  "
$ws.Range("E11").Value = "  -11.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.72"
$ws.Range("E12").Value = "  -7.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("E13").Value = "  -11.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.58"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "4.252.38"
$ws.Range("E15").Value = "  -7.56%  "
$ws.Range("D16").Value = "3.700.65"
$ws.Range("E16").Value = "  -6.70%  "
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("E18").Value = "  -6.00%  "
$ws.Range("E19").Value = "  -9.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.88"
$ws.Range("E20").Value = "  -8.60%  "
$ws.Range("D21").Value = "67.266.37"
$ws.Range("E21").Value = "  -7.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.52"
$ws.Range("E22").Value = "  -7.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.50"
$ws.Range("E23").Value = "  -7.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.59"
$ws.Range("E24").Value = "  -8.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.06"
$ws.Range("E25").Value = "  -9.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.76"
$ws.Range("E26").Value = "  -10.01%  "
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.78"
$ws.Range("E28").Value = "  -7.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.97"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.46"
$ws.Range("E30").Value = "  -9.83%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.72"
$ws.Range("E32").Value = "  -9.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.62"
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("E34").Value = "  -9.96%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.29"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.96"
$ws.Range("E36").Value = "  -7.28%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0910"
$ws.Range("E37").Value = "  -15.03%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "595.59"
$ws.Range("E38").Value = "  -8.19%  "
$ws.Range("E39").Value = "  -7.61%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.17"
$ws.Range("E42").Value = "  +8.58%  "
$ws.Range("E43").Value = "  -7.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.04"
$ws.Range("E44").Value = "  -11.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0436"
$ws.Range("E45").Value = "  -9.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.57"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.35"
$ws.Range("E47").Value = "  -12.80%  "
$ws.Range("E48").Value = "  -9.89%  "
$ws.Range("E49").Value = "  -15.79%  "
$ws.Range("D50").Value = "2.716.41"
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.10"
$ws.Range("E51").Value = "  -9.20%  "
